$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("G9").Value = 0.71
$ws.Range("G10").Select()
